# Modify CFs in aesa_ChangeBiosphereIntegrity_FunctionalDiversity_Hierarchist.xlsx
#
# 1) Highlight only (no value change): "Occupation, forest, intensive" CF
#    (row 14, column C) = 6.1538461538461536E-13 -- leave the higher CF
#    following the precautionary principle.
# 2) Add a new row: "Occupation, lake, artificial" (natural resource::land)
#    = 7.6923076923076923E-13, inserted above the current row 31
#    ("Occupation, shrub land, sclerophyllous"), pushing all following rows
#    down by one. The whole new row is highlighted too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Yellow highlight colour used for both edits (RGB FFFF00 -> BGR 65535)
$yellow = 65535

# --- 1) Highlight the existing "Occupation, forest, intensive" CF cell ---
$ws.Range("C14").Interior.Color = $yellow

# --- 2) Insert the new "Occupation, lake, artificial" row above row 31 ---
$ws.Range("A31").EntireRow.Insert()

$ws.Range("A31").Value2 = "Occupation, lake, artificial"
$ws.Range("B31").Value2 = "natural resource::land"
$ws.Range("C31").Value2 = [double]"7.6923076923076923E-13"

# Highlight the whole new row
$ws.Range("A31:C31").Interior.Color = $yellow
